$d = $word.ActiveDocument

# Update p-values in the Mantel correlogram table (urban 10km)
$d.Content.Find.Execute("0.247", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.248", 2)
$d.Content.Find.Execute("0.494", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.496", 2)
$d.Content.Find.Execute("0.45", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.375", 2)
$d.Content.Find.Execute("0.731", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.713", 2)
